# Update the activity matching report table (B2:L18) with refreshed values.
# This corresponds to the author's fix for stale/unrefreshed chart source data
# (SettingWithCopy related refresh issue) - the underlying numbers for each
# activity row are recomputed/updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(957,13,1400,1504,17,4,932,29202,23,72,1836),
    @(258,0,0,0,0,0,0,0,0,0,0),
    @(4,4,99,109,6,0,550,4242,0,0,0),
    @(1350,1,136,145,7,1,1429,27786,2,0,0),
    @(1181,3,327,351,20,2,3523,48027,4,0,0),
    @(124,0,44,46,2,0,435,7442,1,0,0),
    @(809,1,789,829,5,3,108,16783,18,33,596),
    @(257,0,10,10,0,0,0,8000,0,0,0),
    @(1233,2,627,707,67,10,4584,10597,9,1,33),
    @(1026,1,186,190,3,0,197,15492,0,0,0),
    @(1383,0,89,113,21,3,2059,4336,0,0,0),
    @(860,0,364,364,0,0,0,16157,9,0,0),
    @(259,0,25,25,0,0,0,8000,1,0,0),
    @(364,2,78,80,0,0,0,17700,2,0,0),
    @(181,0,3,3,0,0,0,3333,0,0,0),
    @(1710,1,47,48,0,0,0,16545,4,0,0),
    @(3987,5,229,535,15,1,891,31702,10,283,15833)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 2).Value = $rowValues[$j]
    }
}
